$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.662.10'
$ws.Range('E2').Value = '  -0.03%  '
$ws.Range('D3').Value = '1.597.49'
$ws.Range('E3').Value = '  -0.04%  '
$ws.Range('E4').Value = '  +0.24%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '211.46'
$ws.Range('E5').Value = '  +0.12%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.515'
$ws.Range('E6').Value = '  +0.54%  '
$ws.Range('E7').Value = '  +0.19%  '
$ws.Range('E8').Value = '  +0.09%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.248'
$ws.Range('E9').Value = '  +0.88%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '19.59'
$ws.Range('E10').Value = '  -0.55%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0841'
$ws.Range('E11').Value = '  +0.38%  '
$ws.Range('D12').Value = '1.821.28'
$ws.Range('E12').Value = '  -0.04%  '
$ws.Range('D13').Value = '1.566.38'
$ws.Range('E13').Value = '  -1.97%  '
$ws.Range('E14').Value = '  -0.12%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.524'
$ws.Range('E15').Value = '  +0.42%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '65.05'
$ws.Range('E16').Value = '  +0.36%  '
$ws.Range('D17').Value = '26.646.31'
$ws.Range('E18').Value = '  +1.31%  '
$ws.Range('E19').Value = '  +0.20%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '208.99'
$ws.Range('E20').Value = '  -0.50%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '7.00'
$ws.Range('E21').Value = '  +3.42%  '
$ws.Range('E22').Value = '  +0.45%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '2.32'
$ws.Range('E23').Value = '  +1.12%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '8.99'
$ws.Range('E24').Value = '  +0.93%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '144.23'
$ws.Range('E25').Value = '  -1.36%  '
$ws.Range('E26').Value = '  +0.26%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '7.13'
$ws.Range('E27').Value = '  -0.83%  '
$ws.Range('E28').Value = '  -0.63%  '
$ws.Range('E29').Value = '  -0.19%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.0516'
$ws.Range('E30').Value = '  +2.35%  '
$ws.Range('E31').Value = '  +0.44%  '
$ws.Range('E32').Value = '  +0.25%  '
$ws.Range('E33').Value = '  +1.44%  '
$ws.Range('D34').Value = '1.291.10'
$ws.Range('E34').Value = '  -0.39%  '
$ws.Range('B35').Value = 'ImmutableX'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.615'
$ws.Range('E35').Value = '  -8.12%  '
$ws.Range('B36').Value = 'HuobiToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '2.46'
$ws.Range('E36').Value = '  +0.47%  '
$ws.Range('E37').Value = '  +0.54%  '
$ws.Range('E38').Value = '  -0.78%  '
$ws.Range('E39').Value = '  -0.93%  '
$ws.Range('E40').Value = '  +17.43%  '
$ws.Range('E41').Value = '  +2.18%  '
$ws.Range('E42').Value = '  -0.09%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.783'
$ws.Range('E43').Value = '  -0.47%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '63.48'
$ws.Range('E44').Value = '  -0.57%  '
$ws.Range('D45').Value = '1.733.77'
$ws.Range('E45').Value = '  -0.04%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '90.80'
$ws.Range('E46').Value = '  +0.80%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '1.57'
$ws.Range('E47').Value = '  -3.55%  '
$ws.Range('E48').Value = '  +1.34%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.0508'
$ws.Range('E49').Value = '  +0.83%  '
$ws.Range('E50').Value = '  +0.21%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '7.40'
$ws.Range('E51').Value = '  -1.19%  '
